$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB3").Value = 14.5
$ws.Range("AF3").Value = 19
$ws.Range("AO3").Value = 22
$ws.Range("G3").Value = 2.48
$ws.Range("H3").Value = 2.9
$ws.Range("J3").Value = 3.95
$ws.Range("K3").Value = 4
$ws.Range("P3").Value = 2.42
$ws.Range("W3").Value = 1.67
$ws.Range("X3").Value = 22
$ws.Range("F4").Value = 2.24
$ws.Range("M4").Value = 1.07
$ws.Range("AC5").Value = 7
$ws.Range("V5").Value = 1.71
$ws.Range("AO6").Value = 75
$ws.Range("F7").Value = 6
$ws.Range("R7").Value = 1.56
$ws.Range("S7").Value = 2.68
$ws.Range("V7").Value = 2.64
$ws.Range("F9").Value = 1.09
$ws.Range("G9").Value = 3.1
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 1000
$ws.Range("K9").Value = 950
$ws.Range("S9").Value = 2.92
$ws.Range("V9").Value = 1.01
$ws.Range("W9").Value = 1.47
$ws.Range("G10").Value = 1.36
$ws.Range("O10").Value = 1.16
$ws.Range("S10").Value = 2.22
$ws.Range("F11").Value = 1.17
$ws.Range("I11").Value = 25
$ws.Range("K11").Value = 10.5
$ws.Range("S11").Value = 1.94
$ws.Range("U11").Value = 1.74
$ws.Range("AO12").Value = 36
$ws.Range("G12").Value = 2.8
$ws.Range("P12").Value = 1.88
$ws.Range("T12").Value = 1.76
$ws.Range("W12").Value = 1.55
$ws.Range("G13").Value = 2.08
$ws.Range("K13").Value = 3.5
$ws.Range("W13").Value = 1.92
$ws.Range("AE14").Value = 170
$ws.Range("G14").Value = 1.53
$ws.Range("X14").Value = 17
$ws.Range("F16").Value = 2.68
$ws.Range("I16").Value = 3.3
$ws.Range("V16").Value = 1.43
$ws.Range("W16").Value = 1.55
$ws.Range("P18").Value = 1.7
$ws.Range("G19").Value = 2.36
$ws.Range("H19").Value = 3.55
$ws.Range("W19").Value = 1.73
$ws.Range("H20").Value = 3.45
$ws.Range("K20").Value = 7.8
$ws.Range("AO23").Value = 48
$ws.Range("L23").Value = 1.22
$ws.Range("U23").Value = 1.94
$ws.Range("I24").Value = 2.86
$ws.Range("N25").Value = 2.86
$ws.Range("Q25").Value = 1.86
$ws.Range("AC26").Value = 13
$ws.Range("AF27").Value = 110
$ws.Range("AL27").Value = 110
$ws.Range("AN27").Value = 130
$ws.Range("H27").Value = 1.38
$ws.Range("AA28").Value = 120
$ws.Range("AB28").Value = 8.4
$ws.Range("AD28").Value = 19.5
$ws.Range("AG28").Value = 10.5
$ws.Range("AH28").Value = 20
$ws.Range("AK28").Value = 21
$ws.Range("AN28").Value = 14
$ws.Range("F28").Value = 1.85
$ws.Range("G28").Value = 1.94
$ws.Range("H28").Value = 4.6
$ws.Range("I28").Value = 5.1
$ws.Range("J28").Value = 3.65
$ws.Range("O28").Value = 1.34
$ws.Range("P28").Value = 1.85
$ws.Range("R28").Value = 1.33
$ws.Range("S28").Value = 3.55
$ws.Range("W28").Value = 2.06
$ws.Range("Y28").Value = 16
$ws.Range("Z28").Value = 36
$ws.Range("AD30").Value = 18.5
$ws.Range("G31").Value = 2.38
$ws.Range("H31").Value = 3.5
$ws.Range("J31").Value = 2.9
$ws.Range("L31").Value = 1.38
$ws.Range("N31").Value = 2.76
$ws.Range("O31").Value = 1.39
$ws.Range("P31").Value = 1.7
$ws.Range("Q31").Value = 2
$ws.Range("R31").Value = 1.26
$ws.Range("S31").Value = 3.6
$ws.Range("T31").Value = 1.87
$ws.Range("U31").Value = 1.9
$ws.Range("M33").Value = 1.04
$ws.Range("P33").Value = 1.98
$ws.Range("I34").Value = 4.3
$ws.Range("L34").Value = 1.42
$ws.Range("N34").Value = 2.84
$ws.Range("G35").Value = 2.72
$ws.Range("H35").Value = 2.78
$ws.Range("I35").Value = 3.4
$ws.Range("W35").Value = 1.58
$ws.Range("I37").Value = 14
$ws.Range("T37").Value = 2.14
$ws.Range("AJ38").Value = 42
$ws.Range("AL38").Value = 38
$ws.Range("G38").Value = 2.8
$ws.Range("W38").Value = 1.55
$ws.Range("AM39").Value = 130
$ws.Range("P39").Value = 1.72
$ws.Range("Q39").Value = 2.34
$ws.Range("W39").Value = 1.3
$ws.Range("AM40").Value = 140
$ws.Range("J40").Value = 5.3
$ws.Range("K40").Value = 5.4
$ws.Range("L40").Value = 1.33
$ws.Range("F42").Value = 6
$ws.Range("G42").Value = 7.2
$ws.Range("H42").Value = 1.58
$ws.Range("I42").Value = 1.69
$ws.Range("J42").Value = 3.9
$ws.Range("Q42").Value = 1.87
$ws.Range("V42").Value = 2.42
$ws.Range("N43").Value = 3.3
$ws.Range("J44").Value = 3.25
$ws.Range("W44").Value = 1.41
$ws.Range("F45").Value = 5.2
$ws.Range("G45").Value = 6.2
$ws.Range("H45").Value = 1.67
$ws.Range("I45").Value = 1.81
$ws.Range("L45").Value = 1.32
$ws.Range("Q45").Value = 1.83
$ws.Range("T45").Value = 1.82
$ws.Range("U45").Value = 1.98
$ws.Range("V45").Value = 2.22
$ws.Range("F46").Value = 1.63
$ws.Range("G46").Value = 1.77
$ws.Range("M46").Value = 1.09
$ws.Range("O46").Value = 1.46
$ws.Range("Q46").Value = 2.28
$ws.Range("F47").Value = 2.06
$ws.Range("I47").Value = 3.9
$ws.Range("P47").Value = 2.1
$ws.Range("V47").Value = 1.35
$ws.Range("W47").Value = 1.78
$ws.Range("AF48").Value = 8.800000000000001
$ws.Range("F48").Value = 1.54
$ws.Range("G48").Value = 1.58
$ws.Range("K48").Value = 4.5
$ws.Range("N48").Value = 3.65
$ws.Range("O48").Value = 1.33
$ws.Range("P48").Value = 1.92
$ws.Range("Q48").Value = 1.97
$ws.Range("R48").Value = 1.35
$ws.Range("S48").Value = 3.45
$ws.Range("U48").Value = 1.81
$ws.Range("W48").Value = 2.72
$ws.Range("G49").Value = 3.35
$ws.Range("V49").Value = 1.62
$ws.Range("AB50").Value = 11.5
$ws.Range("AE50").Value = 46
$ws.Range("AI50").Value = 50
$ws.Range("Y50").Value = 19.5
$ws.Range("R51").Value = 1.6
$ws.Range("AE52").Value = 55
$ws.Range("G52").Value = 1.97
$ws.Range("W52").Value = 2.02
$ws.Range("AI53").Value = 120
$ws.Range("N53").Value = 3.6
$ws.Range("P53").Value = 1.9
$ws.Range("R55").Value = 1.35
$ws.Range("G56").Value = 2.44
$ws.Range("J56").Value = 3.4
